# Automatische test-sync: 2025-07-23 14:03:50
# Adds the new "Logs" row (row 4) captured from the 14:03:45 test mail run,
# extends the conditional formatting ranges to cover the new row, and bumps
# the "Dashboard" tally for the "Openingstijden / Locatie" category.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new log entry as row 4 -------------------------------------

$logs.Range("A4").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Range("D4").Value = "Openingstijden / Locatie"
$logs.Range("E4").Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F4").Value = "2025-07-23 14:03:45"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Ja"

# The multi-line content in E4 makes the engine auto-expand the row height;
# restore it to the sheet's standard (non-custom) height like rows 2 and 3.
$logs.Rows.Item(4).EntireRow.AutoFit()

# --- Extend the conditional formatting sqref ranges so row 4 is included ---
# (each block of cfRules sharing a sqref is extended from row 2:3 to 2:4)

foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "3")
    $newRange = $logs.Range($col + "2:" + $col + "4")
    $conditions = $oldRange.FormatConditions
    $count = $conditions.Count()
    for ($i = 1; $i -le $count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard count for "Openingstijden / Locatie" -------------

$dashboard.Range("B2").Value = 3
